$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B78 to be numeric 3 instead of text "3"
$ws.Range("B78").Value = 3

# Add new row 79
$ws.Range("A79").Value = "Ying Tang"
$ws.Range("B79").NumberFormat = "@"
$ws.Range("B79").Value = "4"
$ws.Range("B79").ClearFormats()
$ws.Range("C79").Value = "significant"
$ws.Range("D79").Value = "FBK"
$ws.Range("E79").Value = "OTH"
$ws.Range("F79").Value = "b9d28a3e-28bc-41b5-b6f1-68624390902f"
$ws.Range("G79").Value = "BJlrSmbAZ_annotated.xlsx"
$ws.Range("H79").Value = "The paper's contributions are significant."
